# Ventas.xlsx edit:
#  - Add a new "Date" header column (I1) with the same style as the other headers.
#  - Add a new row of data (row 2) for "Nuevo producto" including a text date.
#  - Update selection / window bookkeeping to match the saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell I1 = "Date" (reuse the bold/bordered header style from H1) ---
$ws.Range("I1").Value = "Date"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats

# --- New data row 2 ---
$ws.Range("A2").Value = 30
$ws.Range("B2").Value = "Nuevo producto"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 12
$ws.Range("F2").Value = 12
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 7
$ws.Range("I2").Value = "22/06/2022"

# --- Window / selection bookkeeping (matches the saved workbook state) ---
$ws.Range("L6").Select()

$win = $wb.Windows.Item(1)
$win.Left = 0
$win.Top = 0
$win.Width = 10410
$win.Height = 1080
